$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calc")

# Insert a new "Latest Morgan" poll row, shifting the old Morgan polls down one.
# Third Morgan (row 8) data is overwritten by what was Second Morgan (row 7).
$secondMorgan = $ws.Range("B7:G7").Value()
$ws.Range("B8:G8").Value = $secondMorgan

# Second Morgan (row 7) data is overwritten by what was Latest Morgan (row 6).
$latestMorgan = $ws.Range("B6:G6").Value()
$ws.Range("B7:G7").Value = $latestMorgan

# Latest Morgan (row 6) gets the brand new poll numbers.
$ws.Range("B6").Value = 55
$ws.Range("C6").Value = 53.5
$ws.Range("D6").Value = 58
$ws.Range("E6").Value = 51.5
$ws.Range("F6").Value = 49
$ws.Range("G6").Value = 58

# Updated Essential 3-poll averages.
$ws.Range("B12").Value = 51.785714285714285
$ws.Range("C12").Value = 50
$ws.Range("D12").Value = 50.179211469534046
$ws.Range("E12").Value = 53.597122302158269
$ws.Range("F12").Value = 52.631578947368418
$ws.Range("G12").Value = 51.578947368421055

$ws.Range("B13").Value = 51.785714285714285
$ws.Range("C13").Value = 49.820788530465954
$ws.Range("D13").Value = 51.971326164874561
$ws.Range("E13").Value = 50.533807829181498
$ws.Range("F13").Value = 53.763440860215056
$ws.Range("G13").Value = 52.173913043478258

$excel.Calculate()

$ws.Range("I19").Select()
